$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.077.25"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.884.19"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'588.29"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'139.58"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.492"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "'6.84"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").Value = "'0.429"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "'0.0000218"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "'32.28"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "3.357.54"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "60.987.15"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "2.885.86"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "'6.50"
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("D19").Value = "'424.87"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "'13.29"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'0.653"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").Value = "'6.91"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'79.76"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "'10.49"
$ws.Range("E24").Value = "  -4.20%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'2.05"
$ws.Range("E26").Value = "  -6.32%  "
$ws.Range("D27").Value = "'11.30"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("D29").Value = "'2.06"
$ws.Range("E29").Value = "  -9.15%  "
$ws.Range("D30").Value = "'6.68"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'25.60"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("D34").Value = "0.0₃0849"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").Value = "'0.967"
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("D36").Value = "'5.43"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("D37").Value = "'48.88"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  -6.31%  "
$ws.Range("D39").Value = "'1.90"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("D40").Value = "'0.116"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("D41").Value = "'8.33"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").Value = "'38.86"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").Value = "'0.263"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").Value = "2.659.43"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "'132.71"
$ws.Range("D46").Value = "'0.0330"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("D47").Value = "'343.62"
$ws.Range("E47").Value = "  -9.57%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "'22.43"
$ws.Range("E49").Value = "  -5.85%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "'1.93"
$ws.Range("E51").Value = "  -3.81%  "
